$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.680.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.39%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.608.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.81%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'565.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.08%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'142.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.81%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.22%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.23%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.632.33"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.34%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'6.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.95%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.67%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.158"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.31%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.374"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +7.90%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.074.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.06%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'60.602.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.19%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'23.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.50%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +1.74%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.620.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.24%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.01%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'10.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +6.94%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'347.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.86%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +12.71%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.14%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +14.13%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'63.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.56%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.11%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.161"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.39%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +6.21%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.0₃0793"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.35%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +4.38%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.11%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'6.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.84%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'160.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.91%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'19.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.78%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'4.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +4.34%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.956"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +8.02%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +4.44%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +5.14%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'37.67"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.74%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.855"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.93%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D42").Value = "'301.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.11%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'141.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +13.61%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.995"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.28%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0983"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.55%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.603"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.16%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0548"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.87%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0241"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +3.76%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +0.65%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'19.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +4.61%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'Maker"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'2.013.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +3.37%  "
$ws.Range("E51").Style = "Normal"
